$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new record row (row 2) with the MCH167-1 collection metadata
$ws.Range("A2").Value = "MCH167-1"
$ws.Range("C2").Value = "ANTI-APARTHEID ACTIVITIES BY MUNICIPALITY AMSTERDAM"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 22C | GRAP COUNT NUMER: NONE"

# Apply the data-row formatting (10pt Calibri, automatic/theme text color) used for
# the new row, matching the rest of the data rows in this export format.
$ws.Range("C2:H2").Font.ThemeColor = 1
$ws.Range("C2:H2").Font.Name = "Calibri"

$ws.Range("A2").Font.ThemeColor = 1
$ws.Range("A2").Font.Name = "Calibri"

# Restore the frozen header pane / selection on the newly active row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A2:J2").Select()
